# data_splitting_plan.pptx — "updating data splitting plan again"
#
# The diagram's group ("Group 4") loses its three right-most "A2"/"A3"/"A4"
# rectangles, and the first rectangle ("Rectangle 5" / the "A1" box) is
# widened (914400 -> 3657600 EMU, i.e. 72pt -> 288pt) so it now spans the
# space the removed boxes used to occupy. Regrouping the remaining shapes
# is what PowerPoint assigns the new group id/name (11 / "Group 10") to.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$grp = $s.Shapes.Item(1)

# Break the group apart so the member shapes become addressable/deletable
# top-level shapes on the slide.
[void]$grp.Ungroup()

# Drop the three trailing rectangles (the "A2", "A3", "A4" boxes).
[void]$s.Shapes.Item("Rectangle 1").Delete()
[void]$s.Shapes.Item("Rectangle 2").Delete()
[void]$s.Shapes.Item("Rectangle 3").Delete()

# Stretch the remaining "A1" rectangle so it now covers the freed-up width.
$rect5 = $s.Shapes.Item("Rectangle 5")
$rect5.Width = 288

# Re-collect every shape left on the slide and regroup them back together.
$names = @()
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $names += $s.Shapes.Item($i).Name
}
$range = $s.Shapes.Range($names)
[void]$range.Group()
